# Generate Report for Handoff
#
# A new source file (dcfdb086-89ff-448d-b135-7d8d584c8f47.md) replaces the
# previously tracked one (c5207598-3f1b-481b-a304-b1da5bbe6d3d.md), a
# dependent file (a2035c89-4d62-4576-a058-71adbf6cabaa.md) failed to
# transform during handoff and is now reported as a new row, and the
# ".localization-config" row shifts down to make room for it.

$wb = $excel.ActiveWorkbook

# Blue hyperlink color (FF6495ED) expressed as a BGR OLE color value.
$hyperlinkColor = 15570276
$dateNumberFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "dcfdb086-89ff-448d-b135-7d8d584c8f47.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("A3").Value = "a2035c89-4d62-4576-a058-71adbf6cabaa.md"
$wsOverview.Range("B3").Value = "Handoff transform failed"
$wsOverview.Range("C3").Value = "Handoff transform failed"

$wsOverview.Range("A4").Value = ".localization-config"
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/e2e/dcfdb086-89ff-448d-b135-7d8d584c8f47.md", "", "", "dcfdb086-89ff-448d-b135-7d8d584c8f47.md")
$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/e2e/a2035c89-4d62-4576-a058-71adbf6cabaa.md", "", "", "a2035c89-4d62-4576-a058-71adbf6cabaa.md")
$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/.localization-config", "", "", ".localization-config")

Style-AsHyperlink $wsOverview.Range("A2:A4")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Range("A2").Value = "dcfdb086-89ff-448d-b135-7d8d584c8f47.md"
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "dcfdb086-89ff-448d-b135-7d8d584c8f47.87dc6c760a6ce5bfd73f390013c0599fcd0372de.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-02-18 03:57:34"
$wsZhCn.Range("G2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Include"

$wsZhCn.Range("A3").Value = "a2035c89-4d62-4576-a058-71adbf6cabaa.md"
$wsZhCn.Range("B3").Value = "Handoff transform failed"
$wsZhCn.Range("D3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("G3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H3").Value = "Ignored"

$wsZhCn.Range("A4").Value = ".localization-config"
$wsZhCn.Range("B4").Value = "Not to be localized"
$wsZhCn.Range("D4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("G4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H4").Value = "Ignored"

$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/e2e/dcfdb086-89ff-448d-b135-7d8d584c8f47.md", "", "", "dcfdb086-89ff-448d-b135-7d8d584c8f47.md")
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e048746b792d0ef07b1be3d47b018996a05163d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/dcfdb086-89ff-448d-b135-7d8d584c8f47.87dc6c760a6ce5bfd73f390013c0599fcd0372de.zh-cn.xlf", "", "", "dcfdb086-89ff-448d-b135-7d8d584c8f47.87dc6c760a6ce5bfd73f390013c0599fcd0372de.zh-cn.xlf")
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/e2e/a2035c89-4d62-4576-a058-71adbf6cabaa.md", "", "", "a2035c89-4d62-4576-a058-71adbf6cabaa.md")
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/.localization-config", "", "", ".localization-config")

Style-AsHyperlink $wsZhCn.Range("A2")
Style-AsHyperlink $wsZhCn.Range("C2")
Style-AsHyperlink $wsZhCn.Range("A3")
Style-AsHyperlink $wsZhCn.Range("A4")

$wsZhCn.Range("D4").NumberFormat = $dateNumberFormat

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Range("A2").Value = "dcfdb086-89ff-448d-b135-7d8d584c8f47.md"
$wsDeDe.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "dcfdb086-89ff-448d-b135-7d8d584c8f47.87dc6c760a6ce5bfd73f390013c0599fcd0372de.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-02-18 03:57:46"
$wsDeDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Include"

$wsDeDe.Range("A3").Value = "a2035c89-4d62-4576-a058-71adbf6cabaa.md"
$wsDeDe.Range("B3").Value = "Handoff transform failed"
$wsDeDe.Range("D3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H3").Value = "Ignored"

$wsDeDe.Range("A4").Value = ".localization-config"
$wsDeDe.Range("B4").Value = "Not to be localized"
$wsDeDe.Range("D4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H4").Value = "Ignored"

$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/e2e/dcfdb086-89ff-448d-b135-7d8d584c8f47.md", "", "", "dcfdb086-89ff-448d-b135-7d8d584c8f47.md")
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6aa582624a4f30721c037eea0df6b750fac2a5c3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/dcfdb086-89ff-448d-b135-7d8d584c8f47.87dc6c760a6ce5bfd73f390013c0599fcd0372de.de-de.xlf", "", "", "dcfdb086-89ff-448d-b135-7d8d584c8f47.87dc6c760a6ce5bfd73f390013c0599fcd0372de.de-de.xlf")
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/e2e/a2035c89-4d62-4576-a058-71adbf6cabaa.md", "", "", "a2035c89-4d62-4576-a058-71adbf6cabaa.md")
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/01380e46bb287ae33a46dddabda686b4aa40b571/.localization-config", "", "", ".localization-config")

Style-AsHyperlink $wsDeDe.Range("A2")
Style-AsHyperlink $wsDeDe.Range("C2")
Style-AsHyperlink $wsDeDe.Range("A3")
Style-AsHyperlink $wsDeDe.Range("A4")

$wsDeDe.Range("D4").NumberFormat = $dateNumberFormat

$wsOverview.Activate()
